# New local strings added: PleaseEnterMessage / PleaseEnterSmsCellPhoneNumbers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 93: PleaseEnterMessage -> Lutfen gonderilecek mesaji giriniz
$ws.Range("A93").Value = "PleaseEnterMessage"
$ws.Range("B93").Value = "Lutfen gonderilecek mesaji giriniz"

# Row 94: PleaseEnterSmsCellPhoneNumbers -> Lutfen mesaj gonderilecek olan telefon numaralarini giriniz.
$ws.Range("A94").Value = "PleaseEnterSmsCellPhoneNumbers"
$ws.Range("B94").Value = "Lutfen mesaj gonderilecek olan telefon numaralarini giriniz."

# Match formatting used by the rest of the key/value table:
# column B (and most of column A) use the style from row 92,
# while A94 follows the alternate style used e.g. by A91.
$ws.Range("B92").Copy()
$ws.Range("A93:B93").PasteSpecial(-4122)
$ws.Range("B94").PasteSpecial(-4122)

$ws.Range("A91").Copy()
$ws.Range("A94").PasteSpecial(-4122)

$excel.CutCopyMode = 0
